# "Made shared car travel count as main shift"
#
# Adds a new "dummy" internal-driver row to the "Internal drivers" sheet
# representing a German shared-car travel slot that should count towards the
# main shift count: Internal driver name = "German dummy 1", Hours per week =
# 0, Country qualifications = "Germany", Home address = "Emmerich,
# Duitsland", Notes = "Dummy".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Internal drivers")
$ws2 = $wb.Worksheets.Item("External driver companies")

# New row 24 on the "Internal drivers" sheet.
$ws1.Cells.Item(24, 1).Value = "German dummy 1"
$ws1.Cells.Item(24, 2).Value = 0
$ws1.Cells.Item(24, 3).Value = "Germany"
$ws1.Cells.Item(24, 4).Value = "Emmerich, Duitsland"
$ws1.Cells.Item(24, 5).Value = "Dummy"

# Match the styling of the rest of the data rows (left/vertical-center align,
# same as every other row in the table, e.g. row 22/23 above it).
$newRow = $ws1.Range("A24:E24")
$newRow.HorizontalAlignment = -4131  # xlHAlignLeft
$newRow.VerticalAlignment = -4108    # xlVAlignCenter

# Restore the cursor positions recorded in the workbook after the edit.
$ws1.Activate()
$ws1.Range("B26").Select()

$ws2.Activate()
$ws2.Range("C24").Select()

$ws1.Activate()
